$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(48, 60, 1.24, 1.25, 5.7, 5.8, 0, 0, 0, 0, 0, 0, 1.51, 2.88, 0, 0, 5.1, 1.01, 1000, 1000, 1000, 1000, 1000, 1000, 1.72, 590, 1000, 1000, 590, 1000, 1000, 1000, 1000, 1000, 1000, 1000)

# Columns F (6) through AO (41) on row 2
$startCol = 6
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(2, $col).Value = $values[$i]
}

$wb.Save()
